# Extracting Text (LEFT, MID, RIGHT)
# Populate columns L (Floor), M (Extension) and N (Wing) by splitting the
# "Location" text held in column K (e.g. "02-West 2635") using the text
# functions LEFT, MID and RIGHT.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First row of data gets its own (non-shared) formula...
$ws.Range("L4").Formula = "=LEFT(K4,2)"
$ws.Range("M4").Formula = "=RIGHT(K4,4)"
$ws.Range("N4").Formula = "=MID(K4,4,4)"

# ...then it is filled down through the rest of the table (rows 5:38),
# which is how the workbook ends up with shared formulas covering L5:L38,
# M5:M38 and N5:N38.
$ws.Range("L5:L38").Formula = "=LEFT(K5,2)"
$ws.Range("M5:M38").Formula = "=RIGHT(K5,4)"
$ws.Range("N5:N38").Formula = "=MID(K5,4,4)"

# Leave the sheet scrolled/selected on the last cell that was filled in,
# matching the state the workbook was left in after the edit.
$ws.Activate()
$ws.Application.ActiveWindow.ScrollColumn = 10
$ws.Range("N14").Select()
